# Move/resize the testimonial callout box on slide 5 and switch its
# run font from "Apple Chancery" to "Bookman Old Style".
#
# Target geometry (EMU):
#   off  x=894520  y=3160642
#   ext cx=9939131 cy=2145268
#
# Shape.Left/Top/Width/Height are expressed in points (1 pt = 12700 EMU);
# the literals below are chosen so the EMU round-trip lands on the exact
# target values.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$shape = $s.Shapes.Item(1)

$shape.Left   = 70.43464660644531
$shape.Top    = 248.8694610595703
$shape.Width  = 782.6087646484375
$shape.Height = 168.9187469482422

$tr = $shape.TextFrame.TextRange
$tr.Font.Name = "Bookman Old Style"
